$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (Tipo de piso / Flat_type.Code),
# shifting the existing columns B:N to D:P.
$ws.Range("B1:C1").EntireColumn.Insert()

# New "Edificio" / Building.Code / Building.Name header block in the
# freshly inserted columns B and C.
$ws.Range("B1").Value = "Edificio"
$ws.Range("B2").Value = "Building.Code"
$ws.Range("C2").Value = "Building.Name"

# Re-create the AutoFilter over the new, wider header range (A2:M3).
$ws.AutoFilterMode = $false
$null = $ws.Range("A2:M3").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$wb.Names.Item(1).RefersTo = "=Precios!`$A`$2:`$M`$3"
